# 966-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment-Makerepayment1.xlsx
# Re-apply "makerepayment" run: a late-repayment recalculation that shifted a
# few pennies of interest/principal across the schedule, plus the UI view
# state (active sheet/tab + selected cell) each sheet was left on.

$wb = $excel.ActiveWorkbook

$wsInput   = $wb.Worksheets.Item("Input")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay   = $wb.Worksheets.Item("Repayment Schedule")
$wsTrans   = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------
# Summary sheet: row 3 (first repayment row) interest/balance corrections
# ---------------------------------------------------------------------
$wsSummary.Range("A3").Value = 672.36
$wsSummary.Range("E3").Value = 570.44000000000005

# ---------------------------------------------------------------------
# Repayment Schedule sheet: recalculated balances for installments 4-12
# ---------------------------------------------------------------------
$wsRepay.Range("F6").Value  = 812.55
$wsRepay.Range("G6").Value  = 6808.56
$wsRepay.Range("H6").Value  = 75.17

$wsRepay.Range("F7").Value  = 818.33
$wsRepay.Range("G7").Value  = 5990.23
$wsRepay.Range("H7").Value  = 69.39

$wsRepay.Range("F8").Value  = 828.64
$wsRepay.Range("G8").Value  = 5161.59
$wsRepay.Range("H8").Value  = 59.08

$wsRepay.Range("G9").Value  = 4326.4799999999996

$wsRepay.Range("F10").Value = 843.63
$wsRepay.Range("G10").Value = 3482.85
$wsRepay.Range("H10").Value = 44.09

$wsRepay.Range("F11").Value = 853.37
$wsRepay.Range("G11").Value = 2629.48
$wsRepay.Range("H11").Value = 34.35

$wsRepay.Range("G12").Value = 1768.56

$wsRepay.Range("F13").Value = 870.28
$wsRepay.Range("G13").Value = 898.28
$wsRepay.Range("H13").Value = 17.440000000000001

$wsRepay.Range("F14").Value = 898.28
$wsRepay.Range("K14").Value = 907.44
$wsRepay.Range("P14").Value = 907.44

# Row 2 had a handful of placeholder cells (A2, E2, N2, O2) cleared out
# entirely in the update.
$wsRepay.Range("A2").ClearContents()
$wsRepay.Range("E2").ClearContents()
$wsRepay.Range("N2").ClearContents()
$wsRepay.Range("O2").ClearContents()

# ---------------------------------------------------------------------
# View state: selection on each sheet, and which sheet/cell is active.
# Order matters -- the last Activate()/Select() wins for the workbook's
# activeTab, and each sheet keeps whatever selection it had when it was
# last active.
# ---------------------------------------------------------------------
$wsInput.Activate()
$wsInput.Range("C19").Select()

$wsSummary.Activate()
$wsSummary.Range("E29").Select()

$wsRepay.Activate()
$wsRepay.Range("I10").Select()

$wsTrans.Activate()
$wsTrans.Range("B3").Select()
